$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "47.251.53"
$ws.Range("E2").Value = "  -1.28%  "
$ws.Range("D3").Value = "2.487.80"
$ws.Range("E3").Value = "  -0.97%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "321.50"
$ws.Range("E5").Value = "  -0.83%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.39"
$ws.Range("E6").Value = "  +2.35%  "
$ws.Range("E7").Value = "  -0.56%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -0.38%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.24"
$ws.Range("E10").Value = "  +3.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0810"
$ws.Range("E11").Value = "  -0.67%  "
$ws.Range("E12").Value = "  +0.71%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.53"
$ws.Range("E13").Value = "  +0.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.17"
$ws.Range("E14").Value = "  -0.16%  "
$ws.Range("D15").Value = "2.875.68"
$ws.Range("E15").Value = "  -0.81%  "
$ws.Range("D16").Value = "2.486.93"
$ws.Range("E16").Value = "  -0.43%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.844"
$ws.Range("E17").Value = "  -0.46%  "
$ws.Range("D18").Value = "47.185.14"
$ws.Range("E18").Value = "  -1.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.46"
$ws.Range("E19").Value = "  +6.04%  "
$ws.Range("E20").Value = "  +0.57%  "
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").Value = "0.0₃0939"
$ws.Range("E21").Value = "  +0.11%  "
$ws.Range("B22").Value = "ImmutableX"
$ws.Range("C22").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.77"
$ws.Range("E22").Value = "  +15.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.47"
$ws.Range("E23").Value = "  -0.48%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "246.29"
$ws.Range("E24").Value = "  -2.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.55"
$ws.Range("E25").Value = "  -0.78%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.69"
$ws.Range("E27").Value = "  -2.47%  "
$ws.Range("E28").Value = "  +4.29%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.95"
$ws.Range("E29").Value = "  -1.00%  "
$ws.Range("E30").Value = "  +2.88%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.59"
$ws.Range("E31").Value = "  -1.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.85"
$ws.Range("E32").Value = "  +0.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.31"
$ws.Range("E33").Value = "  +0.80%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.32"
$ws.Range("E34").Value = "  -1.16%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0781"
$ws.Range("E35").Value = "  -0.28%  "
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.78"
$ws.Range("E37").Value = "  +2.66%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.96"
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("E39").Value = "  -2.25%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "23.09"
$ws.Range("E40").Value = "  +8.82%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.112"
$ws.Range("E41").Value = "  -0.18%  "
$ws.Range("E42").Value = "  -1.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "120.02"
$ws.Range("E43").Value = "  -1.45%  "
$ws.Range("E44").Value = "  -0.46%  "
$ws.Range("D45").Value = "1.998.48"
$ws.Range("E45").Value = "  +1.29%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.02"
$ws.Range("E46").Value = "  +0.77%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.13"
$ws.Range("E48").Value = "  -0.82%  "
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.78"
$ws.Range("E49").Value = "  -2.16%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.12"
$ws.Range("E50").Value = "  -4.89%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "56.83"
$ws.Range("E51").Value = "  +3.25%  "
